# Scheduled market-data refresh: update cached price / profit figures
# on each crafting-job sheet (currentAveragePrice*, Leve* columns H:N).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 357.08
$ws.Range("I15").Value = 357.08
$ws.Range("K15").Value = 1071.24
$ws.Range("M15").Value = -902.24
$ws.Range("H40").Value = 1456.9546
$ws.Range("I40").Value = 1241
$ws.Range("J40").Value = 2191.2
$ws.Range("K40").Value = 1241
$ws.Range("L40").Value = 2191.2
$ws.Range("M40").Value = -1066
$ws.Range("N40").Value = -2541.2
$ws.Range("H96").Value = 557.5263
$ws.Range("I96").Value = 552.38464
$ws.Range("J96").Value = 568.6667
$ws.Range("K96").Value = 1657.15392
$ws.Range("L96").Value = 1706.0001
$ws.Range("M96").Value = -284.15392
$ws.Range("N96").Value = -4452.0001
$ws.Range("H100").Value = 2398.842
$ws.Range("I100").Value = 1250.4166
$ws.Range("J100").Value = 4367.5713
$ws.Range("K100").Value = 1250.4166
$ws.Range("L100").Value = 4367.5713
$ws.Range("M100").Value = -709.4166
$ws.Range("N100").Value = -5449.5713
$ws.Range("H106").Value = 4492.778
$ws.Range("I106").Value = 4595.4
$ws.Range("K106").Value = 4595.4
$ws.Range("M106").Value = -3964.4
$ws.Range("H135").Value = 3234.484
$ws.Range("I135").Value = 2277.3684
$ws.Range("J135").Value = 4749.9165
$ws.Range("K135").Value = 20496.3156
$ws.Range("L135").Value = 42749.2485
$ws.Range("M135").Value = -17961.3156
$ws.Range("N135").Value = -47819.2485

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 364612.03
$ws.Range("I32").Value = 2735.8157
$ws.Range("K32").Value = 2735.8157
$ws.Range("M32").Value = -2448.8157
$ws.Range("H61").Value = 1179.8667
$ws.Range("J61").Value = 1473.5
$ws.Range("L61").Value = 1473.5
$ws.Range("N61").Value = -1897.5
$ws.Range("H132").Value = 1725.2391
$ws.Range("I132").Value = 952.1875
$ws.Range("J132").Value = 3492.2144
$ws.Range("K132").Value = 2856.5625
$ws.Range("L132").Value = 10476.6432
$ws.Range("M132").Value = -326.5625
$ws.Range("N132").Value = -15536.6432
$ws.Range("H136").Value = 1179.8667
$ws.Range("J136").Value = 1473.5
$ws.Range("L136").Value = 4420.5
$ws.Range("N136").Value = -9520.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 67460.13
$ws.Range("I107").Value = 538.53845
$ws.Range("J107").Value = 502450.5
$ws.Range("K107").Value = 538.53845
$ws.Range("L107").Value = 502450.5
$ws.Range("M107").Value = 1381.46155
$ws.Range("N107").Value = -506290.5
$ws.Range("H134").Value = 5858.026
$ws.Range("I134").Value = 1542.5555
$ws.Range("J134").Value = 15567.833
$ws.Range("K134").Value = 4627.666499999999
$ws.Range("L134").Value = 46703.499
$ws.Range("M134").Value = -2092.666499999999
$ws.Range("N134").Value = -51773.499

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2923.38
$ws.Range("I31").Value = 1347.1522
$ws.Range("J31").Value = 21050
$ws.Range("K31").Value = 1347.1522
$ws.Range("L31").Value = 21050
$ws.Range("M31").Value = -1052.1522
$ws.Range("N31").Value = -21640
$ws.Range("H34").Value = 2923.38
$ws.Range("I34").Value = 1347.1522
$ws.Range("J34").Value = 21050
$ws.Range("K34").Value = 1347.1522
$ws.Range("L34").Value = 21050
$ws.Range("M34").Value = -1145.1522
$ws.Range("N34").Value = -21454
$ws.Range("H50").Value = 8666.666999999999
$ws.Range("J50").Value = 8666.666999999999
$ws.Range("L50").Value = 8666.666999999999
$ws.Range("N50").Value = -9916.666999999999
$ws.Range("H104").Value = 10281
$ws.Range("J104").Value = 10281
$ws.Range("L104").Value = 10281
$ws.Range("N104").Value = -15523
$ws.Range("H122").Value = 1561.72
$ws.Range("I122").Value = 1560.1578
$ws.Range("J122").Value = 1566.6666
$ws.Range("K122").Value = 4680.4734
$ws.Range("L122").Value = 4699.9998
$ws.Range("M122").Value = -2230.4734
$ws.Range("N122").Value = -9599.9998
$ws.Range("H132").Value = 1928.7906
$ws.Range("I132").Value = 1695.8975
$ws.Range("J132").Value = 4199.5
$ws.Range("K132").Value = 5087.6925
$ws.Range("L132").Value = 12598.5
$ws.Range("M132").Value = -2557.6925
$ws.Range("N132").Value = -17658.5
$ws.Range("H134").Value = 2021.52
$ws.Range("I134").Value = 1936.4348
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5809.3044
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -3274.3044
$ws.Range("N134").Value = -14070

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 11124444
$ws.Range("J105").Value = 16674166
$ws.Range("L105").Value = 50022498
$ws.Range("N105").Value = -50027740
$ws.Range("H121").Value = 676666.3
$ws.Range("I121").Value = 30000
$ws.Range("J121").Value = 999999.5
$ws.Range("K121").Value = 90000
$ws.Range("L121").Value = 2999998.5
$ws.Range("M121").Value = -88690
$ws.Range("N121").Value = -3002618.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 1000000
$ws.Range("J33").Value = 1000000
$ws.Range("L33").Value = 1000000
$ws.Range("N33").Value = -1000504
$ws.Range("H70").Value = 10334.079
$ws.Range("I70").Value = 13452.941
$ws.Range("K70").Value = 13452.941
$ws.Range("M70").Value = -13182.941
$ws.Range("H73").Value = 10334.079
$ws.Range("I73").Value = 13452.941
$ws.Range("K73").Value = 13452.941
$ws.Range("M73").Value = -12516.941
$ws.Range("H104").Value = 27900
$ws.Range("J104").Value = 27900
$ws.Range("L104").Value = 27900
$ws.Range("N104").Value = -34888
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 1891.238
$ws.Range("I122").Value = 1785.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5357.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2907.4
$ws.Range("N122").Value = -16900

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6900.2104
$ws.Range("I46").Value = 2010.2
$ws.Range("J46").Value = 12333.556
$ws.Range("K46").Value = 2010.2
$ws.Range("L46").Value = 12333.556
$ws.Range("M46").Value = -1822.2
$ws.Range("N46").Value = -12709.556
$ws.Range("H48").Value = 5800
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H55").Value = 822.9259
$ws.Range("I55").Value = 937.8
$ws.Range("J55").Value = 679.3333
$ws.Range("K55").Value = 937.8
$ws.Range("L55").Value = 679.3333
$ws.Range("M55").Value = -764.8
$ws.Range("N55").Value = -1025.3333
$ws.Range("H100").Value = 2930
$ws.Range("J100").Value = 3608.6365
$ws.Range("L100").Value = 3608.6365
$ws.Range("N100").Value = -4690.636500000001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 750.2105
$ws.Range("I107").Value = 489.69232
$ws.Range("J107").Value = 1314.6666
$ws.Range("K107").Value = 1469.07696
$ws.Range("L107").Value = 3943.9998
$ws.Range("M107").Value = 450.9230400000001
$ws.Range("N107").Value = -7783.9998
$ws.Range("H113").Value = 22727678
$ws.Range("I113").Value = 407.9091
$ws.Range("J113").Value = 90909490
$ws.Range("K113").Value = 1223.7273
$ws.Range("L113").Value = 272728470
$ws.Range("M113").Value = 946.2727
$ws.Range("N113").Value = -272732810
$ws.Range("H122").Value = 1182.2759
$ws.Range("I122").Value = 1160.8334
$ws.Range("J122").Value = 1217.3636
$ws.Range("K122").Value = 3482.5002
$ws.Range("L122").Value = 3652.0908
$ws.Range("M122").Value = -1032.5002
$ws.Range("N122").Value = -8552.0908
